$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 and J1 - copy formatting (bold/border/alignment) from H1
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Data values for I2:J31 (I0 and IF columns)
$data = @(
    @(8, 8),
    @(8, 9),
    @(6, 6),
    @(7, 8),
    @(7, 7),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(6, 6),
    @(4, 4),
    @(11, 11),
    @(9, 9),
    @(7, 7),
    @(10, 10),
    @(9, 9),
    @(7, 8),
    @(6, 6),
    @(7, 8),
    @(8, 8),
    @(6, 8),
    @(8, 9),
    @(8, 8),
    @(6, 7),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(5, 7),
    @(1, 2),
    @(3, 3)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
